# NYPD CompStat weekly report update: new crime data collected.
# Updates the report header (volume/number + date range) and the
# Crime Complaints table (rows 15-28, columns C:N) with the latest
# week-to-date / 28-day / year-to-date / 2-year figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Header text updates (rich-text shared strings flatten to plain
# text on write, but the rendered content/formatting is unchanged
# since every run in these strings shares the same font/size/color).
# -----------------------------------------------------------------
# A8: "Volume 32   Number  5" -> "...Number  6"
$ws.Range("A8").Value = "Volume 32   Number  6"
# C9: "Report Covering the Week  1/27/2025  Through  2/2/2025"
#  -> "...2/3/2025  Through  2/9/2025"
$ws.Range("C9").Value = "Report Covering the Week  2/3/2025  Through  2/9/2025"

# -----------------------------------------------------------------
# Cells that change data type/format (text dash placeholder <-> real
# number) need their style copied from a same-style donor cell
# elsewhere on the sheet (rows 14 / 39, which this edit leaves
# untouched) so the number format id matches the target exactly.
#   C14 -> donor for text "0"     (s=13, t=s, v=20)
#   E14 -> donor for text "***.*" (s=13, t=s, v=21)
#   C39 -> donor for integer style (s=15)
#   K39 -> donor for decimal style (s=14)
# -----------------------------------------------------------------

# Row 15
$ws.Range("C39").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("C39").Copy($ws.Range("F15"))
$ws.Range("F15").Value = 1
$ws.Range("I15").Value = 2
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 100

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 12.5
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 10
$ws.Range("K16").Value = 30
$ws.Range("L16").Value = -27.777777777777
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = -78.333333333333

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -52.631578947368
$ws.Range("I17").Value = 11
$ws.Range("J17").Value = 24
$ws.Range("K17").Value = -54.166666666666
$ws.Range("L17").Value = -47.619047619047
$ws.Range("M17").Value = 22.222222222222
$ws.Range("N17").Value = -47.619047619047

# Row 18
$ws.Range("C14").Copy($ws.Range("C18"))
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -45.454545454545
$ws.Range("I18").Value = 8
$ws.Range("J18").Value = 14
$ws.Range("K18").Value = -42.857142857142
$ws.Range("L18").Value = -52.941176470588
$ws.Range("M18").Value = -61.904761904761
$ws.Range("N18").Value = -92.233009708737

# Row 19
$ws.Range("C19").Value = 15
$ws.Range("E19").Value = 36.363636363636
$ws.Range("F19").Value = 55
$ws.Range("G19").Value = 58
$ws.Range("H19").Value = -5.172413793103
$ws.Range("I19").Value = 78
$ws.Range("J19").Value = 89
$ws.Range("K19").Value = -12.359550561797
$ws.Range("L19").Value = 9.859154929577
$ws.Range("M19").Value = 27.868852459016
$ws.Range("N19").Value = -51.851851851851

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 3
$ws.Range("J20").Value = 7
$ws.Range("K20").Value = -57.142857142857
$ws.Range("L20").Value = -50
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = -93.617021276595

# Row 21 (G.L.A. totals, bold style)
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 4.761904761904
$ws.Range("F21").Value = 83
$ws.Range("G21").Value = 99
$ws.Range("H21").Value = -16.161616161616
$ws.Range("I21").Value = 115
$ws.Range("J21").Value = 144
$ws.Range("K21").Value = -20.138888888888
$ws.Range("L21").Value = -13.533834586466
$ws.Range("M21").Value = 7.476635514018
$ws.Range("N21").Value = -70.886075949367

# Row 22 (TOTAL)
$ws.Range("C39").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 2
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -25
$ws.Range("I22").Value = 4
$ws.Range("K22").Value = -50
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 100

# Row 23 (Transit)
$ws.Range("C23").Value = 1
$ws.Range("C39").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 2
$ws.Range("K39").Copy($ws.Range("E23"))
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 25
$ws.Range("I23").Value = 6
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = 50
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 50

# Row 24 (Housing)
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 8.695652173913
$ws.Range("F24").Value = 126
$ws.Range("G24").Value = 99
$ws.Range("H24").Value = 27.272727272727
$ws.Range("I24").Value = 170
$ws.Range("J24").Value = 135
$ws.Range("K24").Value = 25.925925925925
$ws.Range("L24").Value = 51.785714285714
$ws.Range("M24").Value = 38.211382113821

# Row 25 (Petit Larceny)
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = -15
$ws.Range("F25").Value = 99
$ws.Range("G25").Value = 83
$ws.Range("H25").Value = 19.277108433734
$ws.Range("I25").Value = 133
$ws.Range("J25").Value = 113
$ws.Range("K25").Value = 17.699115044247
$ws.Range("L25").Value = 92.753623188405

# Row 26 (Retail Theft)
$ws.Range("C26").Value = 7
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = 7.407407407407
$ws.Range("I26").Value = 37
$ws.Range("J26").Value = 40
$ws.Range("K26").Value = -7.5
$ws.Range("L26").Value = 19.354838709677
$ws.Range("M26").Value = 42.307692307692

# Row 27 (Misd. Assault)
$ws.Range("C39").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("C39").Copy($ws.Range("F27"))
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 2
$ws.Range("K27").Value = 0

# Row 28 (UCR Rape*)
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 3
$ws.Range("K28").Value = -57.142857142857
$ws.Range("L28").Value = -25
